# Update countries & provincias Spain
# Applies the Covid-19 "Pais" dashboard refresh:
#  - bumps the "Datos actualizados..." timestamp string
#  - updates the per-country numeric stat columns (B..H) for the rows
#    whose underlying source numbers changed between the 03:04 and 04:04
#    snapshots

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title / last-updated timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 04:04"

# Estados Unidos (row 4)
$ws.Range("E4").Value = 1019567
$ws.Range("H4").Value = 78615

# Brasil (row 11)
$ws.Range("B11").Value = 146894
$ws.Range("C11").Value = 1002
$ws.Range("E11").Value = 77580
$ws.Range("G11").Value = 25
$ws.Range("H11").Value = 10017

# China (row 14)
$ws.Range("B14").Value = 82887
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 78046
$ws.Range("E14").Value = 208
$ws.Range("F14").Value = 15

# Mexico (row 21)
$ws.Range("D21").Value = 20314
$ws.Range("E21").Value = 8048

# Corea del Sur (row 41)
$ws.Range("B41").Value = 10840
$ws.Range("C41").Value = 18
$ws.Range("D41").Value = 9568
$ws.Range("E41").Value = 1016

# Cuba (row 81)
$ws.Range("B81").Value = 1771
$ws.Range("C81").Value = 86
$ws.Range("D81").Value = 192
$ws.Range("E81").Value = 1472
$ws.Range("F81").Value = 10
$ws.Range("G81").Value = 2
$ws.Range("H81").Value = 107

# Estonia (row 82)
$ws.Range("B82").Value = 1741
$ws.Range("D82").Value = 1078
$ws.Range("E82").Value = 589
$ws.Range("F82").Value = 8
$ws.Range("H82").Value = 74

# Honduras (row 83)
$ws.Range("B83").Value = 1725
$ws.Range("D83").Value = 704
$ws.Range("E83").Value = 965
$ws.Range("F83").Value = 4
$ws.Range("H83").Value = 56

# Nueva Zelanda (row 87)
$ws.Range("B87").Value = 1492
$ws.Range("C87").Value = 2
$ws.Range("D87").Value = 1368
$ws.Range("E87").Value = 103
$ws.Range("F87").Value = 2

# Somalia (row 95)
$ws.Range("B95").Value = 937
$ws.Range("C95").Value = 74
$ws.Range("D95").Value = 130
$ws.Range("E95").Value = 768
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 3
$ws.Range("H95").Value = 39

# Letonia (row 96)
$ws.Range("D96").Value = 106
$ws.Range("E96").Value = 778
$ws.Range("H96").Value = 44

# Kirguistan (row 97)
$ws.Range("B97").Value = 928
$ws.Range("D97").Value = 464
$ws.Range("E97").Value = 446
$ws.Range("F97").Value = 2
$ws.Range("H97").Value = 18

# Republica de Chipre (row 98)
$ws.Range("B98").Value = 906
$ws.Range("D98").Value = 650
$ws.Range("E98").Value = 244
$ws.Range("F98").Value = 13
$ws.Range("H98").Value = 12

# Consejo Danes para los Refugiados (row 99)
$ws.Range("B99").Value = 891
$ws.Range("D99").Value = 400
$ws.Range("E99").Value = 476
$ws.Range("F99").Value = 10
$ws.Range("H99").Value = 15

# Reunion (row 127)
$ws.Range("B127").Value = 426
$ws.Range("E127").Value = 72

# Guyana (row 164)
$ws.Range("D164").Value = 35
$ws.Range("E164").Value = 49

# Nueva Caledonia (row 192)
$ws.Range("D192").Value = 16
$ws.Range("H192").Value = 2

# Belice (row 193)
$ws.Range("D193").Value = 18
$ws.Range("H193").Value = 0
